$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the test/demo patient data that was filled into the admission form,
# resetting it back to a blank template now that it is in production.

# 1er. Apellido / 2do. Apellido / 1er. Nombre / 2do. Nombre (row 6)
$ws.Range("A6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("G6").ClearContents()

# Dirección actual: Departamento / Teléfono (row 8)
$ws.Range("H8").ClearContents()
$ws.Range("J8").ClearContents()

# Dirección habitual: Calle o lugar / Municipio / Departamento / Teléfono (row 10)
$ws.Range("A10").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("F10").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("J10").ClearContents()

# Fecha de nacimiento / Edad en años / Lugar de nacimiento (row 12)
$ws.Range("A12").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("H12").ClearContents()

# Estado Civil: Soltero -> Casado (row 14)
$ws.Range("A14").Value = "Casado"

# Ocupación / Nacionalidad / No. De Cédula (row 14)
$ws.Range("D14").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("H14").ClearContents()

# Nombre del Padre / Nombre de la Madre (row 18)
$ws.Range("A18").ClearContents()
$ws.Range("F18").ClearContents()

# En caso de emergencia notificar a / Parentesco / Teléfono (row 20)
$ws.Range("A20").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("J20").ClearContents()

# Fecha de Ingreso / Hora (row 24)
$ws.Range("A24").ClearContents()
$ws.Range("C24").ClearContents()
